# Add new translated-indicator rows to the dictionary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New code/description pairs (appended after the current last row, 108).
$newRows = @(
    @{ Code = "CASH_INKASS";      Desc = "Сумма инкассации" },
    @{ Code = "CARD_PURCHASES";   Desc = "Сумма покупок по карте" },
    @{ Code = "WALLET_POTENTIAL"; Desc = "Средний потенциал кошелька" },
    @{ Code = "ESCAPE_RKO";       Desc = "Доля закрывших РКО" }
)

$startRow = 109
$lastRow = $startRow + $newRows.Count - 1

# A column keeps the same bordered/left-aligned look used by the rest of
# the code column, so carry that formatting over from an existing cell.
$formatSource = $ws.Range("A2")

for ($i = 0; $i -lt ($newRows.Count - 1); $i++) {
    $row = $startRow + $i
    $codeCell = $ws.Cells.Item($row, 1)
    $descCell = $ws.Cells.Item($row, 2)

    $formatSource.Copy($codeCell)
    $codeCell.Value = $newRows[$i].Code
    $descCell.Value = $newRows[$i].Desc
}

# The very last row's code cell instead gets a plain bordered style (no
# fill/alignment carried over) to close out the table.
$lastCodeCell = $ws.Cells.Item($lastRow, 1)
$lastDescCell = $ws.Cells.Item($lastRow, 2)
$lastCodeCell.Value = $newRows[$newRows.Count - 1].Code
$lastDescCell.Value = $newRows[$newRows.Count - 1].Desc
$lastCodeCell.Borders.LineStyle = 1

# Select the last newly-added cell, mirroring the saved selection state.
$ws.Range("A$lastRow").Select()
